# Insert a new weekly price record as the new row 427 for the
# "Hortaliza, Feria Lagunitas de Puerto Montt - Repollo" sheet.
# This pushes the previous rows 427:492 down to 428:493 (dimension
# grows from A1:R492 to A1:R493), matching the commit "Fruta / hortaliza,
# semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 427; Excel shifts 427:492 -> 428:493
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new weekly record
$ws.Cells.Item(427, 1).Value  = 4
$ws.Cells.Item(427, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(427, 3).Value  = "Los Lagos"
$ws.Cells.Item(427, 4).Value  = 44776
$ws.Cells.Item(427, 5).Value  = 10
$ws.Cells.Item(427, 6).Value  = 100112006
$ws.Cells.Item(427, 7).Value  = "Repollo"
$ws.Cells.Item(427, 8).Value  = "Crespo record"
$ws.Cells.Item(427, 9).Value  = "Primera"
$ws.Cells.Item(427, 10).Value = 100
$ws.Cells.Item(427, 11).Value = 2000
$ws.Cells.Item(427, 12).Value = 2000
$ws.Cells.Item(427, 13).Value = 2000
$ws.Cells.Item(427, 14).Value = "`$/unidad"
$ws.Cells.Item(427, 15).Value = "Región Metropolitana"
$ws.Cells.Item(427, 16).Value = 2000
$ws.Cells.Item(427, 17).Value = 1
$ws.Cells.Item(427, 18).Value = "Hortaliza"
